$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 480-485: change date from 44209 to 44939 (and some price updates) ---

# Row 480: date + K,L,M,P (1000 -> 1200)
$ws.Range("D480").Value = 44939
$ws.Range("K480").Value = 1200
$ws.Range("L480").Value = 1200
$ws.Range("M480").Value = 1200
$ws.Range("P480").Value = 1200

# Row 481: date + J (6000 -> 4000) + K,L,M,P (800 -> 900)
$ws.Range("D481").Value = 44939
$ws.Range("J481").Value = 4000
$ws.Range("K481").Value = 900
$ws.Range("L481").Value = 900
$ws.Range("M481").Value = 900
$ws.Range("P481").Value = 900

# Row 482: date only
$ws.Range("D482").Value = 44939

# Row 483: date + K,L,M,P (1000 -> 1300)
$ws.Range("D483").Value = 44939
$ws.Range("K483").Value = 1300
$ws.Range("L483").Value = 1300
$ws.Range("M483").Value = 1300
$ws.Range("P483").Value = 1300

# Row 484: date + J (6000 -> 5000) + K,L,M,P (800 -> 900)
$ws.Range("D484").Value = 44939
$ws.Range("J484").Value = 5000
$ws.Range("K484").Value = 900
$ws.Range("L484").Value = 900
$ws.Range("M484").Value = 900
$ws.Range("P484").Value = 900

# Row 485: date only
$ws.Range("D485").Value = 44939

# --- Append new rows 486-491: copies of the ORIGINAL (pre-edit) rows 480-485 ---

$newRows = @(
    @{ Row = 486; H = "Calameño"; I = "Extra";    J = 3000; K = 1000; L = 1000; M = 1000; P = 1000 },
    @{ Row = 487; H = "Calameño"; I = "Primera";  J = 6000; K = 800;  L = 800;  M = 800;  P = 800  },
    @{ Row = 488; H = "Calameño"; I = "Segunda";  J = 4000; K = 600;  L = 600;  M = 600;  P = 600  },
    @{ Row = 489; H = "Tuna";     I = "Extra";    J = 3000; K = 1000; L = 1000; M = 1000; P = 1000 },
    @{ Row = 490; H = "Tuna";     I = "Primera";  J = 6000; K = 800;  L = 800;  M = 800;  P = 800  },
    @{ Row = 491; H = "Tuna";     I = "Segunda";  J = 4000; K = 600;  L = 600;  M = 600;  P = 600  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 5
    $ws.Range("B$row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$row").Value = "Maule"
    $ws.Range("D$row").Value = 44209
    $ws.Range("D$row").NumberFormat = $ws.Range("D480").NumberFormat
    $ws.Range("E$row").Value = 7
    $ws.Range("F$row").Value = 100112027
    $ws.Range("G$row").Value = "Melón"
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = "`$/unidad"
    $ws.Range("O$row").Value = "Región del Maule"
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = 1
    $ws.Range("R$row").Value = "Hortaliza"
}
